# EPBDS-12729: Remove the extra trailing args from the error() formula
# samples and drop the redundant "error3" example block (the one that
# showed "= error('foo.bar', null)" with no trailing args, which is now
# redundant with the trimmed error2 sample). What used to be the
# "error4" example becomes the new "error3".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# error1 block (rows 4-6): drop the trailing ", 1, 2, 3" args
$ws.Range("D6").Value = "'" + '= error("foo.bar", "Foo bar")'

# error2 block (rows 10-12): drop the trailing ", 1, 2, 3" args
$ws.Range("D12").Value = "'" + '= error("foo.bar", null)'

# Remove the old error3 block (rows 15-17) entirely; deleting these rows
# shifts the old error4 block (rows 23-25) up to rows 20-22.
$ws.Rows("15:17").Delete()

# Rename the relocated former-error4 header to error3.
$ws.Range("C20").Value = "Spreadsheet SpreadsheetResult error3()"

# Match the author's final selection position.
$ws.Range("D41").Select()
